# Apply updated crypto price/volume figures (and reordered coin rows) per the commit diff.
# Values are entered with a leading apostrophe to force text interpretation (avoids numeric
# auto-conversion/rounding for cells like "0.487"), then ClearFormats() strips the resulting
# "Text" number-format style so the cell keeps the original unstyled appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.ClearFormats()
}

Set-TextValue "D2" "26.893.67"
Set-TextValue "E2" "  +0.05%  "
Set-TextValue "D3" "1.546.60"
Set-TextValue "E3" "  -1.11%  "
Set-TextValue "E4" "  +0.31%  "
Set-TextValue "D5" "206.62"
Set-TextValue "E5" "  +0.27%  "
Set-TextValue "D6" "0.487"
Set-TextValue "E6" "  -0.22%  "
Set-TextValue "E8" "  -0.22%  "
Set-TextValue "D9" "21.41"
Set-TextValue "E9" "  -1.58%  "
Set-TextValue "D10" "0.0583"
Set-TextValue "E10" "  -0.15%  "
Set-TextValue "D11" "0.0856"
Set-TextValue "E11" "  -1.09%  "
Set-TextValue "D12" "1.766.97"
Set-TextValue "E12" "  -1.01%  "
Set-TextValue "D13" "1.553.32"
Set-TextValue "E13" "  -0.53%  "
Set-TextValue "E14" "  -0.86%  "
Set-TextValue "D15" "0.512"
Set-TextValue "E15" "  -0.41%  "
Set-TextValue "D16" "26.897.48"
Set-TextValue "E16" "  +0.14%  "
Set-TextValue "D17" "61.40"
Set-TextValue "E17" "  +0.36%  "
Set-TextValue "D18" "214.27"
Set-TextValue "E18" "  +0.01%  "
Set-TextValue "D19" "0.0₃0683"
Set-TextValue "E19" "  +0.69%  "
Set-TextValue "E20" "  -1.96%  "
Set-TextValue "E21" "  +0.29%  "
Set-TextValue "E22" "  -2.82%  "
Set-TextValue "E23" "  -0.24%  "
Set-TextValue "E24" "  -3.13%  "
Set-TextValue "D25" "151.70"
Set-TextValue "E25" "  -1.49%  "
Set-TextValue "D26" "6.62"
Set-TextValue "E26" "  -1.23%  "
Set-TextValue "D27" "14.86"
Set-TextValue "E27" "  -0.54%  "
Set-TextValue "E28" "  +0.29%  "
Set-TextValue "E29" "  +0.34%  "
Set-TextValue "E30" "  -0.89%  "
Set-TextValue "E31" "  -0.93%  "
Set-TextValue "E32" "  +2.14%  "
Set-TextValue "D33" "1.365.29"
Set-TextValue "E33" "  -2.61%  "
Set-TextValue "E34" "  +1.32%  "
Set-TextValue "D35" "1.53"
Set-TextValue "E35" "  +0.27%  "
Set-TextValue "E36" "  +4.28%  "
Set-TextValue "E37" "  +0.39%  "
Set-TextValue "E38" "  -0.13%  "
Set-TextValue "D39" "0.520"
Set-TextValue "E39" "  -0.69%  "
Set-TextValue "E40" "  -1.20%  "
Set-TextValue "E41" "  +0.31%  "
Set-TextValue "D42" "5.61"
Set-TextValue "E42" "  +5.30%  "
Set-TextValue "E43" "  -0.73%  "
Set-TextValue "E44" "  +1.96%  "
Set-TextValue "D45" "63.46"
Set-TextValue "E45" "  +0.52%  "
Set-TextValue "E46" "  -2.55%  "
Set-TextValue "B47" "mCoin"
Set-TextValue "C47" "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
Set-TextValue "D47" "2.26"
Set-TextValue "E47" "  -3.62%  "
Set-TextValue "B48" "RocketPoolETH"
Set-TextValue "C48" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D48" "1.681.76"
Set-TextValue "E48" "  -0.90%  "
Set-TextValue "B49" "Quant"
Set-TextValue "C49" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D49" "85.58"
Set-TextValue "E49" "  -0.65%  "
Set-TextValue "B50" "Cronos"
Set-TextValue "C50" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D50" "0.0508"
Set-TextValue "E50" "  +0.61%  "
Set-TextValue "B51" "BabyDogeCoin"
Set-TextValue "C51" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D51" "0.0₇0972"
Set-TextValue "E51" "  -1.21%  "
